$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates (existing row, revised figures) ---
$ws.Range("B2").Value = "'2"
$ws.Range("D2").Value = 0.241
$ws.Range("E2").Value = 0.281
$ws.Range("G2").Value = 0.2922182757507265
$ws.Range("H2").Value = 0.2690776019804111
$ws.Range("I2").Value = 0.2335593585189969
$ws.Range("J2").Value = 0.207747022517156
$ws.Range("K2").Value = 22.33
$ws.Range("L2").Value = 0.2403401140889032
$ws.Range("M2").Value = 6.89
$ws.Range("N2").Value = 0.02917019475021168
$ws.Range("O2").Value = 0.3085535154500672
$ws.Range("P2").Value = 6.89
$ws.Range("Q2").Value = 0.02917019475021168
$ws.Range("R2").Value = 0.3085535154500672
$ws.Range("U2").Value = 15.74
$ws.Range("V2").Value = 0.06663844199830651
$ws.Range("W2").Value = 0.1562621493331994
$ws.Range("X2").Value = 0.03463223085152103
$ws.Range("Y2").Value = 0.1216299184816783
$ws.Range("Z2").Value = 0.4252563163676308
$ws.Range("AA2").Value = 0.08107996797840525
$ws.Range("AB2").Value = 0.02809267379038723
$ws.Range("AC2").Value = 0.05298729418801802
$ws.Range("AD2").Value = 118.8
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 118.8
$ws.Range("AG2").Value = 103.06
$ws.Range("AH2").Value = 0.3346478873239437
$ws.Range("AI2").Value = 0.4177215189873418
$ws.Range("AJ2").Value = 0.3037788127100159
$ws.Range("AK2").Value = 0.3836075336856994
$ws.Range("AL2").Value = 1.87
$ws.Range("AM2").Value = 1.87
$ws.Range("AN2").Value = 5.210526315789474
$ws.Range("AO2").Value = 11.60427807486631
$ws.Range("AP2").Value = 4.520175438596492
$ws.Range("AQ2").Value = 11.60427807486631

# --- Row 3 updates (existing row, revised figures) ---
$ws.Range("D3").Value = 0.241
$ws.Range("E3").Value = 0.281
$ws.Range("G3").Value = 0.3220640569395018
$ws.Range("H3").Value = 0.2965599051008304
$ws.Range("I3").Value = 0.2574139976275208
$ws.Range("J3").Value = 0.2005166278070929
$ws.Range("K3").Value = 19.9
$ws.Range("L3").Value = 0.236061684460261
$ws.Range("M3").Value = 6.89
$ws.Range("N3").Value = 0.03282515483563601
$ws.Range("O3").Value = 0.3462311557788945
$ws.Range("P3").Value = 6.89
$ws.Range("Q3").Value = 0.03282515483563601
$ws.Range("R3").Value = 0.3462311557788945
$ws.Range("U3").Value = 12.8
$ws.Range("V3").Value = 0.06098141972367795
$ws.Range("W3").Value = 0.144938091769847
$ws.Range("X3").Value = 0.02065062324325839
$ws.Range("Y3").Value = 0.1242874685265886
$ws.Range("Z3").Value = 0.808710667689946
$ws.Range("AA3").Value = 0.1621599359568105
$ws.Range("AB3").Value = 0.02092922455977079
$ws.Range("AC3").Value = 0.1412307113970397
$ws.Range("AD3").Value = 16.4
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 16.4
$ws.Range("AG3").Value = 3.599999999999998
$ws.Range("AH3").Value = 0.07247017233760494
$ws.Range("AI3").Value = 0.0984984984984985
$ws.Range("AJ3").Value = 0.01686182669789226
$ws.Range("AK3").Value = 0.02342225113858164
$ws.Range("AL3").Value = 1.87
$ws.Range("AM3").Value = 1.87
$ws.Range("AN3").Value = 0.719298245614035
$ws.Range("AO3").Value = 11.60427807486631
$ws.Range("AP3").Value = 0.1578947368421052
$ws.Range("AQ3").Value = 11.60427807486631

# --- Row 4: new company row appended to the table ---
$ws.Range("A4").Value = "Colombia"
$ws.Range("B4").Value = "Credifamilia Compañía de Financiamiento S.A. (BVC:CREDIFAMI)"
$ws.Range("C4").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2.43
$ws.Range("L4").Value = 0.2822299651567944
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 2.94
$ws.Range("V4").Value = 0.111787072243346
$ws.Range("W4").Value = 0.1675862068965517
$ws.Range("X4").Value = 0.04861383845978366
$ws.Range("Y4").Value = 0.1189723684367681
$ws.Range("Z4").Value = 0.07536764705882353
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.03525612302100366
$ws.Range("AC4").Value = -0.03525612302100366
$ws.Range("AD4").Value = 102.4
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 102.4
$ws.Range("AG4").Value = 99.46000000000001
$ws.Range("AH4").Value = 0.7956487956487955
$ws.Range("AI4").Value = 0.8685326547921968
$ws.Range("AJ4").Value = 0.7908715012722647
$ws.Range("AK4").Value = 0.8651704940848991
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0

